$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.081.20"
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").Value = "1.551.78"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "287.44"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "0.3817"
$ws.Range("E7").Value = "  +2.33%  "
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("D9").Value = "43.58"
$ws.Range("E9").Value = "  -9.68%  "
$ws.Range("D10").Value = "1.125"
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("D11").Value = "0.07346"
$ws.Range("E11").Value = "  -1.56%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "20.05"
$ws.Range("E13").Value = "  -2.82%  "
$ws.Range("D14").Value = "5.769"
$ws.Range("E14").Value = "  -2.76%  "
$ws.Range("D15").Value = "6.739"
$ws.Range("E15").Value = "  -2.40%  "
$ws.Range("D16").Value = "1.561.82"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").Value = "0.00001078"
$ws.Range("E17").Value = "  -3.06%  "
$ws.Range("D18").Value = "0.06638"
$ws.Range("E18").Value = "  -1.70%  "
$ws.Range("D19").Value = "85.64"
$ws.Range("E19").Value = "  -2.31%  "
$ws.Range("D21").Value = "6.342"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "16.00"
$ws.Range("E22").Value = "  -2.49%  "
$ws.Range("D23").Value = "11.65"
$ws.Range("E23").Value = "  -3.41%  "
$ws.Range("D24").Value = "22.068.32"
$ws.Range("E24").Value = "  -1.59%  "
$ws.Range("D25").Value = "2.301"
$ws.Range("E25").Value = "  -3.44%  "
$ws.Range("D26").Value = "2.494"
$ws.Range("E26").Value = "  -2.72%  "
$ws.Range("D27").Value = "150.51"
$ws.Range("E27").Value = "  -1.55%  "
$ws.Range("D28").Value = "19.13"
$ws.Range("E28").Value = "  -2.59%  "
$ws.Range("D29").Value = "4.939"
$ws.Range("E29").Value = "  -1.44%  "
$ws.Range("D30").Value = "121.61"
$ws.Range("E30").Value = "  -1.91%  "
$ws.Range("D31").Value = "1.735.06"
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").Value = "1.074"
$ws.Range("E32").Value = "  +2.12%  "
$ws.Range("D33").Value = "5.855"
$ws.Range("E33").Value = "  -4.25%  "
$ws.Range("D34").Value = "1.901"
$ws.Range("E34").Value = "  -5.44%  "
$ws.Range("D35").Value = "0.08204"
$ws.Range("E35").Value = "  -1.26%  "
$ws.Range("D36").Value = "9.264"
$ws.Range("E36").Value = "  -3.76%  "
$ws.Range("D37").Value = "0.06272"
$ws.Range("E37").Value = "  -1.62%  "
$ws.Range("D38").Value = "0.02309"
$ws.Range("E38").Value = "  -5.93%  "
$ws.Range("D39").Value = "5.247"
$ws.Range("E39").Value = "  -1.98%  "
$ws.Range("D40").Value = "0.2145"
$ws.Range("E40").Value = "  -5.62%  "
$ws.Range("E41").Value = "  -4.39%  "
$ws.Range("D42").Value = "10.99"
$ws.Range("E42").Value = "  -2.17%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "0.6000"
$ws.Range("E44").Value = "  -4.21%  "
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("D46").Value = "3.728"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("D47").Value = "0.5815"
$ws.Range("E47").Value = "  -4.69%  "
$ws.Range("D48").Value = "1.965"
$ws.Range("E48").Value = "  -3.66%  "
$ws.Range("D49").Value = "121.69"
$ws.Range("E49").Value = "  -2.73%  "
$ws.Range("D50").Value = "1.171"
$ws.Range("E50").Value = "  -3.13%  "
$ws.Range("D51").Value = "0.07007"
$ws.Range("E51").Value = "  -2.99%  "
